$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'60.964.23"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.09%  "
$ws.Range("D3").Value = "'2.916.97"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.23%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'592.84"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.14%  "
$ws.Range("D6").Value = "'146.03"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.28%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "'0.507"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.80%  "
$ws.Range("D9").Value = "'6.90"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.48%  "
$ws.Range("D10").Value = "'0.144"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.28%  "
$ws.Range("E11").Value = "  -1.82%  "
$ws.Range("E12").Value = "  +0.55%  "
$ws.Range("D13").Value = "'33.56"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.20%  "
$ws.Range("E14").Value = "  -0.33%  "
$ws.Range("D15").Value = "'3.398.91"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.17%  "
$ws.Range("D16").Value = "'60.934.31"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.11%  "
$ws.Range("E17").Value = "  -1.10%  "
$ws.Range("D18").Value = "'2.916.01"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.18%  "
$ws.Range("D19").Value = "'430.78"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.02%  "
$ws.Range("E20").Value = "  -1.49%  "
$ws.Range("E21").Value = "  +1.09%  "
$ws.Range("E22").Value = "  -0.65%  "
$ws.Range("D23").Value = "'81.60"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.80%  "
$ws.Range("D24").Value = "'10.96"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.11%  "
$ws.Range("D25").Value = "'2.20"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.02%  "
$ws.Range("D26").Value = "'11.91"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.63%  "
$ws.Range("E27").Value = "  +0.03%  "
$ws.Range("D28").Value = "'2.28"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.64%  "
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("E30").Value = "  -0.41%  "
$ws.Range("D31").Value = "'7.04"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.71%  "
$ws.Range("D32").Value = "'26.54"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.32%  "
$ws.Range("E33").Value = "  +1.36%  "
$ws.Range("D34").Value = "'0.0₃0854"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.19%  "
$ws.Range("E35").Value = "  +0.20%  "
$ws.Range("E36").Value = "  +0.15%  "
$ws.Range("E37").Value = "  +1.44%  "
$ws.Range("E38").Value = "  -0.38%  "
$ws.Range("E39").Value = "  -1.46%  "
$ws.Range("D40").Value = "'8.54"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.30%  "
$ws.Range("E41").Value = "  -1.55%  "
$ws.Range("D42").Value = "'39.95"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.84%  "
$ws.Range("D43").Value = "'375.29"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.49%  "
$ws.Range("D44").Value = "'0.0345"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.70%  "
$ws.Range("D45").Value = "'2.705.12"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.29%  "
$ws.Range("D46").Value = "'132.40"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.28%  "
$ws.Range("D48").Value = "'23.85"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.58%  "
$ws.Range("E49").Value = "  -0.24%  "
$ws.Range("E50").Value = "  -3.33%  "
